{"js": "// Update the date line and each two-digit \u00f7 one-digit division problem in\n// the practice-sheet table. Matches are looked up by their exact current\n// text and replaced in place (so cell formatting/runs are preserved).\n// \"80\u00f78=\" appears twice with two different replacements, so that one is\n// handled by resolving both search hits (returned in document order) and\n// replacing them by index.\n\nconst replacements = [\n  [\"2024-03-14 Thursday\", \"2024-03-15 Friday\"],\n  [\"99\u00f78=\", \"58\u00f73=\"],\n  [\"96\u00f79=\", \"73\u00f78=\"],\n  [\"27\u00f73=\", \"89\u00f79=\"],\n  [\"63\u00f75=\", \"71\u00f76=\"],\n  [\"85\u00f74=\", \"92\u00f77=\"],\n  [\"85\u00f79=\", \"59\u00f75=\"],\n  [\"17\u00f75=\", \"40\u00f79=\"],\n  [\"88\u00f72=\", \"21\u00f76=\"],\n  [\"36\u00f74=\", \"48\u00f76=\"],\n  [\"94\u00f76=\", \"46\u00f77=\"],\n  [\"38\u00f79=\", \"90\u00f72=\"],\n  [\"87\u00f78=\", \"51\u00f74=\"],\n  [\"84\u00f75=\", \"41\u00f75=\"],\n  [\"28\u00f75=\", \"77\u00f73=\"],\n  [\"24\u00f77=\", \"26\u00f76=\"],\n  [\"61\u00f73=\", \"61\u00f75=\"],\n  [\"86\u00f73=\", \"29\u00f74=\"],\n  [\"74\u00f78=\", \"30\u00f74=\"],\n  [\"43\u00f77=\", \"60\u00f73=\"],\n  [\"85\u00f77=\", \"71\u00f74=\"],\n  [\"32\u00f72=\", \"78\u00f77=\"],\n  [\"42\u00f76=\", \"60\u00f75=\"],\n  [\"59\u00f73=\", \"77\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"80\u00f78=\" occurs twice (first instance -> \"92\u00f78=\", second -> \"80\u00f74=\").\n// Resolve both hits in one search call; search results come back in\n// document order, so index 0 is the first occurrence and index 1 the\n// second.\nconst dup = context.document.body.search(\"80\u00f78=\", { matchCase: true });\ndup.load(\"items\");\nawait context.sync();\ndup.items[0].insertText(\"92\u00f78=\", Word.InsertLocation.replace);\ndup.items[1].insertText(\"80\u00f74=\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the date line and each two-digit \u00f7 one-digit division problem in\n# the practice-sheet table using Find/Replace against the document\n# content range. Each Find.Execute call only replaces the first remaining\n# match (Replace = 1 / wdReplaceOne), so a fresh Find against\n# $d.Content always targets the next not-yet-updated occurrence in\n# document order \u2014 this is what lets the two \"80\u00f78=\" cells (which need\n# two different replacement values) be handled safely without touching\n# the wrong one.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-14 Thursday\", \"2024-03-15 Friday\"),\n    @(\"99\u00f78=\", \"58\u00f73=\"),\n    @(\"96\u00f79=\", \"73\u00f78=\"),\n    @(\"27\u00f73=\", \"89\u00f79=\"),\n    @(\"80\u00f78=\", \"92\u00f78=\"),\n    @(\"63\u00f75=\", \"71\u00f76=\"),\n    @(\"85\u00f74=\", \"92\u00f77=\"),\n    @(\"85\u00f79=\", \"59\u00f75=\"),\n    @(\"17\u00f75=\", \"40\u00f79=\"),\n    @(\"88\u00f72=\", \"21\u00f76=\"),\n    @(\"36\u00f74=\", \"48\u00f76=\"),\n    @(\"94\u00f76=\", \"46\u00f77=\"),\n    @(\"38\u00f79=\", \"90\u00f72=\"),\n    @(\"80\u00f78=\", \"80\u00f74=\"),\n    @(\"87\u00f78=\", \"51\u00f74=\"),\n    @(\"84\u00f75=\", \"41\u00f75=\"),\n    @(\"28\u00f75=\", \"77\u00f73=\"),\n    @(\"24\u00f77=\", \"26\u00f76=\"),\n    @(\"61\u00f73=\", \"61\u00f75=\"),\n    @(\"86\u00f73=\", \"29\u00f74=\"),\n    @(\"74\u00f78=\", \"30\u00f74=\"),\n    @(\"43\u00f77=\", \"60\u00f73=\"),\n    @(\"85\u00f77=\", \"71\u00f74=\"),\n    @(\"32\u00f72=\", \"78\u00f77=\"),\n    @(\"42\u00f76=\", \"60\u00f75=\"),\n    @(\"59\u00f73=\", \"77\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n}\n\nWrite-Output \"done\"\n"}
